$wb = $excel.ActiveWorkbook

# --- Sheet 1: "dimenzija za m=4" (columns A:B) ---
$ws1 = $wb.Worksheets.Item(1)

# Extend formatting (style/borders/font) of the last existing data row (72)
# down into the three new rows (73:75), matching how the existing rows are
# formatted (bold+bordered column A, plain column B).
$ws1.Range("A72:B72").Copy()
$ws1.Range("A73:B75").PasteSpecial(-4122)

$ws1.Range("A73").Value = 78
$ws1.Range("B73").Value = 32
$ws1.Range("A74").Value = 79
$ws1.Range("B74").Value = 32
$ws1.Range("A75").Value = 80
$ws1.Range("B75").Value = 32

# --- Sheet 2: "dimenzija za m=5-10" (columns A:G) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A72:G72").Copy()
$ws2.Range("A73:G75").PasteSpecial(-4122)

$ws2.Range("A73").Value = 78
$ws2.Range("B73").Value = 26
$ws2.Range("C73").Value = 30
$ws2.Range("D73").Value = 30
$ws2.Range("E73").Value = 26
$ws2.Range("F73").Value = 28
$ws2.Range("G73").Value = 29

$ws2.Range("A74").Value = 79
$ws2.Range("B74").Value = 28
$ws2.Range("C74").Value = 30
$ws2.Range("D74").Value = 30
$ws2.Range("E74").Value = 29
$ws2.Range("F74").Value = 30
$ws2.Range("G74").Value = 29

$ws2.Range("A75").Value = 80
$ws2.Range("B75").Value = 28
$ws2.Range("C75").Value = 30
$ws2.Range("D75").Value = 30
$ws2.Range("E75").Value = 29
$ws2.Range("F75").Value = 30
$ws2.Range("G75").Value = 30

$excel.CutCopyMode = $false
